$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Timestamp" (col A) and "Notified Production (MW)" (col B) values for rows 2-97.
$rowData = @{
    2 = @(46048.01041666666, 1101.716)
    3 = @(46048.02083333334, 1102.039)
    4 = @(46048.03125, 1105.753)
    5 = @(46048.04166666666, 1107.947)
    6 = @(46048.05208333334, 1086.88)
    7 = @(46048.0625, 1083.883)
    8 = @(46048.07291666666, 1080.828)
    9 = @(46048.08333333334, 1090.565)
    10 = @(46048.09375, 1112.264)
    11 = @(46048.10416666666, 1111.065)
    12 = @(46048.11458333334, 1094.499)
    13 = @(46048.125, 1100.095)
    14 = @(46048.13541666666, 1125.717)
    15 = @(46048.14583333334, 1141.02)
    16 = @(46048.15625, 1153.939)
    17 = @(46048.16666666666, 1164.699)
    18 = @(46048.17708333334, 1202.892)
    19 = @(46048.1875, 1223.042)
    20 = @(46048.19791666666, 1243.582)
    21 = @(46048.20833333334, 1269.881)
    22 = @(46048.21875, 1321.709)
    23 = @(46048.22916666666, 1331.1)
    24 = @(46048.23958333334, 1364.465)
    25 = @(46048.25, 1387.2)
    26 = @(46048.26041666666, 1414.398)
    27 = @(46048.27083333334, 1451.946)
    28 = @(46048.28125, 1481.248)
    29 = @(46048.29166666666, 1514.838)
    30 = @(46048.30208333334, 1542.983)
    31 = @(46048.3125, 1569.846)
    32 = @(46048.32291666666, 1609.337)
    33 = @(46048.33333333334, 1652.398)
    34 = @(46048.34375, 1669.671)
    35 = @(46048.35416666666, 1682.026)
    36 = @(46048.36458333334, 1692.539)
    37 = @(46048.375, 1702.331)
    38 = @(46048.38541666666, 1754.731)
    39 = @(46048.39583333334, 1758.813)
    40 = @(46048.40625, 1761.748)
    41 = @(46048.41666666666, 1765.151)
    42 = @(46048.42708333334, 1748.232)
    43 = @(46048.4375, 1751.486)
    44 = @(46048.44791666666, 1754.289)
    45 = @(46048.45833333334, 1757.821)
    46 = @(46048.46875, 1788.058)
    47 = @(46048.47916666666, 1795.429)
    48 = @(46048.48958333334, 1803.093)
    49 = @(46048.5, 1808.888)
    50 = @(46048.51041666666, 1823.814)
    51 = @(46048.52083333334, 1835.011)
    52 = @(46048.53125, 1846.892)
    53 = @(46048.54166666666, 1861.2)
    54 = @(46048.55208333334, 1883.215)
    55 = @(46048.5625, 1895.295)
    56 = @(46048.57291666666, 1906.769)
    57 = @(46048.58333333334, 1918.476)
    58 = @(46048.59375, 1950.239)
    59 = @(46048.60416666666, 1959.19)
    60 = @(46048.61458333334, 1966.8)
    61 = @(46048.625, 1975.237)
    62 = @(46048.63541666666, 1993.795)
    63 = @(46048.64583333334, 2001.035)
    64 = @(46048.65625, 2008.698)
    65 = @(46048.66666666666, 2015.384)
    66 = @(46048.67708333334, 2027.816)
    67 = @(46048.6875, 2037.76)
    68 = @(46048.69791666666, 2046.713)
    69 = @(46048.70833333334, 2054.723)
    70 = @(46048.71875, 2069.757)
    71 = @(46048.72916666666, 2076.09)
    72 = @(46048.73958333334, 2078.482)
    73 = @(46048.75, 2081.26)
    74 = @(46048.76041666666, 2091.872)
    75 = @(46048.77083333334, 2092.214)
    76 = @(46048.78125, 2090.851)
    77 = @(46048.79166666666, 2091.227)
    78 = @(46048.80208333334, 2102.51)
    79 = @(46048.8125, 2105.951)
    80 = @(46048.82291666666, 2109.178)
    81 = @(46048.83333333334, 2114.225)
    82 = @(46048.84375, 2116.677)
    83 = @(46048.85416666666, 2125.973)
    84 = @(46048.86458333334, 2134.939)
    85 = @(46048.875, 2142.235)
    86 = @(46048.88541666666, 2142.705)
    87 = @(46048.89583333334, 2140.639)
    88 = @(46048.90625, 2134.807)
    89 = @(46048.91666666666, 2131.16)
    90 = @(46048.92708333334, 2124.294)
    91 = @(46048.9375, 2114.384)
    92 = @(46048.94791666666, 2104.221)
    93 = @(46048.95833333334, 2094.773)
    94 = @(46048.96875, 0)
    95 = @(46048.97916666666, 0)
    96 = @(46048.98958333334, 0)
    97 = @(46049, 0)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 1).Value2 = $vals[0]
    $ws.Cells.Item($row, 2).Value2 = $vals[1]
}